$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
}

$ws.Range("D2").Value = "42.461.83"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "2.288.36"
$ws.Range("E3").Value = "  +0.55%  "

Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "156.13"
$ws.Range("E5").Value = "  +15,497.08%  "

$ws.Range("D6").Value = "306.27"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "96.63"
$ws.Range("E7").Value = "  +4.38%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("D11").Value = "35.51"
$ws.Range("E11").Value = "  +8.72%  "

$ws.Range("D12").Value = "0.0812"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").Value = "2.643.03"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "2.287.52"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "0.797"
$ws.Range("E18").Value = "  +4.24%  "

$ws.Range("D19").Value = "42.311.19"
$ws.Range("E19").Value = "  +1.24%  "

Set-TextCell "D20" "12.90"
$ws.Range("E20").Value = "  +4.95%  "

$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").Value = "68.23"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").Value = "244.54"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "1.96"
$ws.Range("E26").Value = "  +1.64%  "

Set-TextCell "D28" "24.30"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "36.54"
$ws.Range("E29").Value = "  +6.99%  "

$ws.Range("D30").Value = "9.72"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("E31").Value = "  +1.51%  "

Set-TextCell "D32" "161.30"
$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").Value = "0.0755"
$ws.Range("E35").Value = "  +0.61%  "

Set-TextCell "D36" "3.10"
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "17.48"
$ws.Range("E37").Value = "  +2.13%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.109"
$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  +0.98%  "

$ws.Range("D42").Value = "4.23"
$ws.Range("E42").Value = "  +7.71%  "

$ws.Range("D43").Value = "19.86"
$ws.Range("E43").Value = "  +2.21%  "

$ws.Range("D44").Value = "2.013.18"
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("E45").Value = "  +10.98%  "

$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").Value = "10.28"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("D49").Value = "54.14"
$ws.Range("E49").Value = "  +4.40%  "

$ws.Range("D50").Value = "1.55"
$ws.Range("E50").Value = "  +1.56%  "

$ws.Range("D51").Value = "73.12"
$ws.Range("E51").Value = "  +0.20%  "

